# Auto-generated edit script: updates crypto price/volume table
# to match the scraped data snapshot from Wed Sep 11 23:13:00 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.472.38"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.316.44"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.52"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.40"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("D9").Value = "2.339.14"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.36"
$ws.Range("E12").Value = "  -2.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").Value = "2.755.12"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.34"
$ws.Range("E15").Value = "  -4.33%  "
$ws.Range("D16").Value = "57.348.31"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "2.329.95"
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "335.94"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.40"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.78"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.21"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.72"
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.54"
$ws.Range("E29").Value = "  +4.13%  "
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").Value = "0.0₃0724"
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.11"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.47"
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("E36").Value = "  -3.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.929"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.98"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.16"
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("E40").Value = "  -2.03%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "148.88"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.375"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.60"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "283.47"
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.08"
$ws.Range("E45").Value = "  -3.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0931"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0500"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.81"
$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.558"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0217"
$ws.Range("E50").Value = "  -1.26%  "
$ws.Range("E51").Value = "  +6.57%  "
